$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -5
    5  = 8
    10 = 5
    16 = -2
    18 = -4
    19 = 8
    22 = -2
    26 = -2
    33 = 5
    36 = -3
    38 = 0
    40 = 5
    43 = -2
    44 = -4
    46 = 3
    54 = -6
    55 = 8
    58 = -1
    59 = 3
    61 = 0
    65 = 2
    68 = 3
    70 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
